$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently follows the
#    H1 title paragraph near the top of the document.
# ---------------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Replace the final paragraph (the "Create a feature image..." image
#    prompt paragraph) with two new paragraphs: a bold title line followed
#    by an italic meta-description line.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r/>
<w:r><w:rPr><w:b/></w:rPr><w:t>Play 2 Gods Zeus versus Thor Free | Innovative Dual Spin Mechanism</w:t></w:r>
</w:p>
<w:p>
<w:r/>
<w:r><w:rPr><w:i/></w:rPr><w:t>Experience high volatility and significant rewards with 2 Gods Zeus versus Thor, featuring innovative Dual Spin and Win Spins features. Play for free now!</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$rng.InsertXML($xml)
